$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.050185454231615
$ws.Range("D2").Value = 1.059304086085447
$ws.Range("E2").Value = 1.063983032370308
$ws.Range("F2").Value = 1.071705009219449
$ws.Range("I2").Value = 1.051694525032248
$ws.Range("J2").Value = 1.055219944866521
$ws.Range("K2").Value = 1.062033955755731
$ws.Range("L2").Value = 1.066700198581982
$ws.Range("M2").Value = 1.074401472400465
$ws.Range("N2").Value = 1.022225663660216

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.051046424322753
$ws.Range("D3").Value = 1.060002788512794
$ws.Range("E3").Value = 1.064799386103898
$ws.Range("F3").Value = 1.072530654749289
$ws.Range("I3").Value = 1.051948241710204
$ws.Range("J3").Value = 1.055730928409041
$ws.Range("K3").Value = 1.062547153733213
$ws.Range("L3").Value = 1.067331665881884
$ws.Range("M3").Value = 1.075043699842557
$ws.Range("N3").Value = 1.022397244018378

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.051604179609631
$ws.Range("D4").Value = 1.060455425350396
$ws.Range("E4").Value = 1.065328571255784
$ws.Range("F4").Value = 1.073065807134554
$ws.Range("I4").Value = 1.052111458659547
$ws.Range("J4").Value = 1.056061538132145
$ws.Range("K4").Value = 1.062879078314866
$ws.Range("L4").Value = 1.067740561570622
$ws.Range("M4").Value = 1.075459513837651
$ws.Range("N4").Value = 1.022508214081469

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.051838813974156
$ws.Range("D5").Value = 1.060645838927361
$ws.Range("E5").Value = 1.065551266215851
$ws.Range("F5").Value = 1.073290999860103
$ws.Range("I5").Value = 1.052179845691697
$ws.Range("J5").Value = 1.056200517839244
$ws.Range("K5").Value = 1.063018582454606
$ws.Range("L5").Value = 1.067912530231413
$ws.Range("M5").Value = 1.075634380314305
$ws.Range("N5").Value = 1.022554852520607

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.051878219103477
$ws.Range("D6").Value = 1.060677817486972
$ws.Range("E6").Value = 1.06558867085086
$ws.Range("F6").Value = 1.073328823242868
$ws.Range("I6").Value = 1.052191314701106
$ws.Range("J6").Value = 1.056223852601841
$ws.Range("K6").Value = 1.06304200360798
$ws.Range("L6").Value = 1.067941408523767
$ws.Range("M6").Value = 1.075663744537224
$ws.Range("N6").Value = 1.022562682521717

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.051607314202193
$ws.Range("D7").Value = 1.060457969176764
$ws.Range("E7").Value = 1.065331546032898
$ws.Range("F7").Value = 1.073068815329657
$ws.Range("I7").Value = 1.052112373352713
$ws.Range("J7").Value = 1.056063395221229
$ws.Range("K7").Value = 1.062880942522874
$ws.Range("L7").Value = 1.067742859155383
$ws.Range("M7").Value = 1.075461850185044
$ws.Range("N7").Value = 1.022508837319723

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.05047628755026
$ws.Range("D8").Value = 1.059540105124689
$ws.Range("E8").Value = 1.064258725190701
$ws.Range("F8").Value = 1.071983851574751
$ws.Range("I8").Value = 1.051780466959497
$ws.Range("J8").Value = 1.055392639773672
$ws.Range("K8").Value = 1.062207423433667
$ws.Range("L8").Value = 1.066913544110994
$ws.Range("M8").Value = 1.074618463495442
$ws.Range("N8").Value = 1.022283660923095

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.048488323002812
$ws.Range("D9").Value = 1.057926845112654
$ws.Range("E9").Value = 1.062375632605405
$ws.Range("F9").Value = 1.070079017510505
$ws.Range("I9").Value = 1.051188329648528
$ws.Range("J9").Value = 1.054210502280555
$ws.Range("K9").Value = 1.061019515513299
$ws.Range("L9").Value = 1.065454503987391
$ws.Range("M9").Value = 1.07313429094592
$ws.Range("N9").Value = 1.021886479344543

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.047166502252414
$ws.Range("D10").Value = 1.056854221900368
$ws.Range("E10").Value = 1.061125292049586
$ws.Range("F10").Value = 1.068813953849679
$ws.Range("I10").Value = 1.050788728142699
$ws.Range("J10").Value = 1.053422363435551
$ws.Range("K10").Value = 1.060226928882644
$ws.Range("L10").Value = 1.064483459070362
$ws.Range("M10").Value = 1.072146268048102
$ws.Range("N10").Value = 1.021621456132763

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.046594986605184
$ws.Range("D11").Value = 1.05639047010966
$ws.Range("E11").Value = 1.060585102754073
$ws.Range("F11").Value = 1.068267335566175
$ws.Range("I11").Value = 1.050614557715479
$ws.Range("J11").Value = 1.053081095063205
$ws.Range("K11").Value = 1.059883593043903
$ws.Range("L11").Value = 1.064063395733729
$ws.Range("M11").Value = 1.0717188008225
$ws.Range("N11").Value = 1.021506648213588

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.046382828249971
$ws.Range("D12").Value = 1.056218319159801
$ws.Range("E12").Value = 1.060384637019528
$ws.Range("F12").Value = 1.06806447378175
$ws.Range("I12").Value = 1.05054969254404
$ws.Range("J12").Value = 1.05295433418981
$ws.Range("K12").Value = 1.059756043034561
$ws.Range("L12").Value = 1.063907427946633
$ws.Range("M12").Value = 1.071560075095401
$ws.Range("N12").Value = 1.021463996234793

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.046428331119498
$ws.Range("D13").Value = 1.056255241247307
$ws.Range("E13").Value = 1.060427629215559
$ws.Range("F13").Value = 1.068107980301632
$ws.Range("I13").Value = 1.050563614046969
$ws.Range("J13").Value = 1.052981524744205
$ws.Range("K13").Value = 1.059783403818783
$ws.Range("L13").Value = 1.063940880716573
$ws.Range("M13").Value = 1.071594119805217
$ws.Range("N13").Value = 1.021473145547493

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.046577446908969
$ws.Range("D14").Value = 1.056376237863089
$ws.Range("E14").Value = 1.060568528416884
$ws.Range("F14").Value = 1.068250563334468
$ws.Range("I14").Value = 1.050609199414073
$ws.Range("J14").Value = 1.053070616934161
$ws.Range("K14").Value = 1.059873050116386
$ws.Range("L14").Value = 1.064050502105286
$ws.Range("M14").Value = 1.071705679380808
$ws.Range("N14").Value = 1.021503122731968

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.046669339077349
$ws.Range("D15").Value = 1.056450802118168
$ws.Range("E15").Value = 1.060655365585152
$ws.Range("F15").Value = 1.068338436903152
$ws.Range("I15").Value = 1.050637263491783
$ws.Range("J15").Value = 1.053125509795147
$ws.Range("K15").Value = 1.059928281576268
$ws.Range("L15").Value = 1.06411805177998
$ws.Range("M15").Value = 1.071774422203627
$ws.Range("N15").Value = 1.021521591724044

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.047204449727319
$ws.Range("D16").Value = 1.056885014503164
$ws.Range("E16").Value = 1.061161168422089
$ws.Range("F16").Value = 1.068850255763422
$ws.Range("I16").Value = 1.05080026330322
$ws.Range("J16").Value = 1.053445012420367
$ws.Range("K16").Value = 1.060249712102917
$ws.Range("L16").Value = 1.064511345963692
$ws.Range("M16").Value = 1.072174645215563
$ws.Range("N16").Value = 1.021629074519711

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.047540336899748
$ws.Range("D17").Value = 1.057157573429053
$ws.Range("E17").Value = 1.061478772167165
$ws.Range("F17").Value = 1.069171618831429
$ws.Range("I17").Value = 1.050902204043824
$ws.Range("J17").Value = 1.05364542911393
$ws.Range("K17").Value = 1.060451300277952
$ws.Range("L17").Value = 1.064758158652837
$ws.Range("M17").Value = 1.07242579020795
$ws.Range("N17").Value = 1.021696482286521

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.04773633528703
$ws.Range("D18").Value = 1.057316619969874
$ws.Range("E18").Value = 1.061664142299438
$ws.Range("F18").Value = 1.069359176520656
$ws.Range("I18").Value = 1.050961554273032
$ws.Range("J18").Value = 1.053762328787972
$ws.Range("K18").Value = 1.060568869594921
$ws.Range("L18").Value = 1.064902159343917
$ws.Range("M18").Value = 1.072572312807371
$ws.Range("N18").Value = 1.02173579513235

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.047803179376911
$ws.Range("D19").Value = 1.057370862130498
$ws.Range("E19").Value = 1.061727368553425
$ws.Range("F19").Value = 1.069423147792976
$ws.Range("I19").Value = 1.050981772471552
$ws.Range("J19").Value = 1.053802188519033
$ws.Range("K19").Value = 1.060608955337135
$ws.Range("L19").Value = 1.064951266431546
$ws.Range("M19").Value = 1.072622278942243
$ws.Range("N19").Value = 1.021749198940411

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.047504290986795
$ws.Range("D20").Value = 1.057128323437948
$ws.Range("E20").Value = 1.061444684150796
$ws.Range("F20").Value = 1.069137128012061
$ws.Range("I20").Value = 1.050891278147421
$ws.Range("J20").Value = 1.053623926289623
$ws.Range("K20").Value = 1.060429673177368
$ws.Range("L20").Value = 1.06473167393961
$ws.Range("M20").Value = 1.072398841201634
$ws.Range("N20").Value = 1.02168925058112

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.046533532479344
$ws.Range("D21").Value = 1.056340604404383
$ws.Range("E21").Value = 1.060527032002442
$ws.Range("F21").Value = 1.06820857129514
$ws.Range("I21").Value = 1.050595780360366
$ws.Range("J21").Value = 1.053044381456209
$ws.Range("K21").Value = 1.059846652051772
$ws.Range("L21").Value = 1.064018219602751
$ws.Range("M21").Value = 1.071672826345601
$ws.Range("N21").Value = 1.021494295393408

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.045923919448167
$ws.Range("D22").Value = 1.055845954788792
$ws.Range("E22").Value = 1.059951136764651
$ws.Range("F22").Value = 1.067625773441298
$ws.Range("I22").Value = 1.050409002997449
$ws.Range("J22").Value = 1.052680007309469
$ws.Range("K22").Value = 1.059479969680303
$ws.Range("L22").Value = 1.063570004549704
$ws.Range("M22").Value = 1.071216668858935
$ws.Range("N22").Value = 1.021371677873365

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.046247016047176
$ws.Range("D23").Value = 1.056108118490282
$ws.Range("E23").Value = 1.060256327765179
$ws.Range("F23").Value = 1.067934628040598
$ws.Range("I23").Value = 1.050508110420785
$ws.Range("J23").Value = 1.052873167696459
$ws.Range("K23").Value = 1.059674365290191
$ws.Range("L23").Value = 1.063807577010039
$ws.Range("M23").Value = 1.071458456014291
$ws.Range("N23").Value = 1.021436683511917

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.047520578332933
$ws.Range("D24").Value = 1.057141540041442
$ws.Range("E24").Value = 1.06146008669607
$ws.Range("F24").Value = 1.069152712582721
$ws.Range("I24").Value = 1.050896215429928
$ws.Range("J24").Value = 1.053633642490476
$ws.Range("K24").Value = 1.060439445575547
$ws.Range("L24").Value = 1.0647436411211
$ws.Range("M24").Value = 1.072411018192799
$ws.Range("N24").Value = 1.021692518292902

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.049001651459211
$ws.Range("D25").Value = 1.058343410889397
$ws.Range("E25").Value = 1.062861574305481
$ws.Range("F25").Value = 1.070570620417796
$ws.Range("I25").Value = 1.051342268565253
$ws.Range("J25").Value = 1.054516126917853
$ws.Range("K25").Value = 1.061326738184251
$ws.Range("L25").Value = 1.065831416818277
$ws.Range("M25").Value = 1.073517740651626
$ws.Range("N25").Value = 1.021989203864477
